{"js": "// Map of old equation text -> new equation text, per the commit diff.\nconst replacements = {\n  \"59\u00d724=\": \"50\u00d749=\",\n  \"43\u00d796=\": \"96\u00d768=\",\n  \"62\u00d769=\": \"30\u00d798=\",\n  \"52\u00d773=\": \"61\u00d785=\",\n  \"28\u00d773=\": \"50\u00d712=\",\n  \"37\u00d783=\": \"36\u00d757=\",\n  \"56\u00d783=\": \"66\u00d790=\",\n  \"81\u00d747=\": \"74\u00d796=\",\n  \"17\u00d757=\": \"16\u00d770=\",\n  \"19\u00d799=\": \"39\u00d736=\",\n  \"89\u00d753=\": \"88\u00d754=\",\n  \"87\u00d759=\": \"50\u00d743=\",\n  \"77\u00d745=\": \"37\u00d738=\",\n  \"64\u00d758=\": \"76\u00d761=\",\n  \"34\u00d764=\": \"31\u00d796=\",\n  \"74\u00d762=\": \"73\u00d712=\",\n  \"63\u00d780=\": \"95\u00d789=\",\n  \"84\u00d730=\": \"28\u00d790=\",\n  \"31\u00d737=\": \"52\u00d774=\",\n  \"99\u00d717=\": \"97\u00d765=\",\n  \"81\u00d756=\": \"62\u00d731=\",\n  \"64\u00d775=\": \"64\u00d787=\",\n  \"17\u00d726=\": \"80\u00d742=\",\n  \"50\u00d798=\": \"34\u00d716=\",\n  \"87\u00d756=\": \"96\u00d769=\",\n};\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nfor (const p of paragraphs.items) {\n  const text = p.text;\n  if (Object.prototype.hasOwnProperty.call(replacements, text)) {\n    p.insertText(replacements[text], Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Old equation text -> new equation text, per the commit diff.\n$replacements = [ordered]@{\n    \"59\u00d724=\" = \"50\u00d749=\"\n    \"43\u00d796=\" = \"96\u00d768=\"\n    \"62\u00d769=\" = \"30\u00d798=\"\n    \"52\u00d773=\" = \"61\u00d785=\"\n    \"28\u00d773=\" = \"50\u00d712=\"\n    \"37\u00d783=\" = \"36\u00d757=\"\n    \"56\u00d783=\" = \"66\u00d790=\"\n    \"81\u00d747=\" = \"74\u00d796=\"\n    \"17\u00d757=\" = \"16\u00d770=\"\n    \"19\u00d799=\" = \"39\u00d736=\"\n    \"89\u00d753=\" = \"88\u00d754=\"\n    \"87\u00d759=\" = \"50\u00d743=\"\n    \"77\u00d745=\" = \"37\u00d738=\"\n    \"64\u00d758=\" = \"76\u00d761=\"\n    \"34\u00d764=\" = \"31\u00d796=\"\n    \"74\u00d762=\" = \"73\u00d712=\"\n    \"63\u00d780=\" = \"95\u00d789=\"\n    \"84\u00d730=\" = \"28\u00d790=\"\n    \"31\u00d737=\" = \"52\u00d774=\"\n    \"99\u00d717=\" = \"97\u00d765=\"\n    \"81\u00d756=\" = \"62\u00d731=\"\n    \"64\u00d775=\" = \"64\u00d787=\"\n    \"17\u00d726=\" = \"80\u00d742=\"\n    \"50\u00d798=\" = \"34\u00d716=\"\n    \"87\u00d756=\" = \"96\u00d769=\"\n}\n\nforeach ($old in $replacements.Keys) {\n    $new = $replacements[$old]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.Text = $new\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n\n    $find.Execute([ref]$old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null\n}\n"}
